# ydf_treeple_scaling.xlsx edit script
# Implements: "Add sparse_oblique_exponent to YDF. Add more d values to MIGHT
# experiments. Update XLSX"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Top "YDF" table (rows 11-17): add "0 repeats" label at E11, and a new
#    "m5.metal" note above it at B9.
# ---------------------------------------------------------------------------
$ws.Range("B9").Value = "m5.metal"
$ws.Range("E11").Value = "0 repeats"

# ---------------------------------------------------------------------------
# 2) Third table (rows 21-27): this used to be the "Treeple" table with only
#    a "Kernel keeps restarting" note. It becomes a second YDF run
#    ("double-check work of YDF" / "YDF - sparse_oblique_exponent = 1"),
#    now fully populated with data, and an "0 repeats" label.
# ---------------------------------------------------------------------------
$ws.Range("I19").Value = "double-check work of YDF"
$ws.Range("A20").Value = "YDF - sparse_oblique_exponent = 1"

# A21 flips from "Treeple" to "YDF"
$ws.Range("A21").Value = "YDF"
$ws.Range("E21").Value = "0 repeats"

# B23 used to hold the text "Kernel keeps restarting" - replace with real data
$ws.Range("B23").Value = 0.3935

$tbl3 = @{
  "C23" = 0.968;    "D23" = 2.7609;   "E23" = 6.0874;   "F23" = 20.6349;  "G23" = 73.3892;  "H23" = 203.8382
  "B24" = 0.8763;   "C24" = 2.1466;   "D24" = 5.9235;   "E24" = 12.6541;  "F24" = 41.807;   "G24" = 146.6751; "H24" = 399.8987
  "B25" = 2.0816;   "C25" = 4.9963;   "D25" = 12.973;   "E25" = 26.9844;  "F25" = 85.4866;  "G25" = 295.1661; "H25" = 802.4954
  "B26" = 140.0341; "C26" = 166.26;   "D26" = 222.1643; "E26" = 296.4236; "F26" = 480.165
}
foreach ($addr in $tbl3.Keys) {
  $ws.Range($addr).Value = $tbl3[$addr]
}

# ---------------------------------------------------------------------------
# 3) New fourth table (rows 30-36): Treeple, now with d values up to 131072
#    (160*2^10) via doubling formulas, and "3 repeats" label.
# ---------------------------------------------------------------------------
$ws.Range("A30").Value = "Treeple"
$ws.Range("B30").Value = "d"
$ws.Range("C30").Value = "keeping n_attributes = 160"
$ws.Range("E30").Value = "3 repeats"

$ws.Range("A31").Value = "n"
$ws.Range("B31").Value = 160
$ws.Range("C31").Value = 320
$ws.Range("D31").Value = 640
$ws.Range("E31").Value = 1024
$ws.Range("F31").Value = 2048
$ws.Range("G31").Value = 4096
$ws.Range("H31").Value = 8192

# Make the header row (B31:H31) match the bold-italic style used by the
# other table headers (A2:H2, A12:H12, A22:H22, ...).
$ws.Range("B31:H31").Font.Bold = $true
$ws.Range("B31:H31").Font.Italic = $true

# Doubling formulas out to d=131072, continuing the n_attributes progression.
$ws.Range("I31").Formula = "=8192*2"
$ws.Range("J31").Formula = "=I31*2"
$ws.Range("K31:L31").Formula = "=J31*2"

# I31 gets its own bold/italic (no explicit font colour) style, distinct from
# the rest of the header row.
$ws.Range("I31").Font.Bold = $true
$ws.Range("I31").Font.Italic = $true
$ws.Range("I31").Font.ColorIndex = 0

# J31:L31 should share the same style as B31:H31 (bold italic, themed colour).
$ws.Range("J31:L31").Font.Bold = $true
$ws.Range("J31:L31").Font.Italic = $true

$tbl4 = @{
  "A32" = 500;  "B32" = 7.1137;  "C32" = 7.2495;  "D32" = 7.4532;  "E32" = 7.5957;  "F32" = 8.0981;  "G32" = 10.3894; "H32" = 20.5901
  "A33" = 1000; "B33" = 10.1264; "C33" = 10.262;  "D33" = 10.4505; "E33" = 10.7079; "F33" = 11.5363; "G33" = 16.9053; "H33" = 38.5697
  "A34" = 2000; "B34" = 15.675;  "C34" = 15.718;  "D34" = 16.1071; "E34" = 16.4595; "F34" = 18.5542; "G34" = 32.3944; "H34" = 76.3515
  "A35" = 4000; "B35" = 25.874;  "C35" = 26.1485; "D35" = 26.8861; "E35" = 27.7313; "F35" = 31.8937; "G35" = 64.9342; "H35" = 160.3889
  "A36" = 8000; "B36" = 45.749;  "C36" = 46.4161; "D36" = 48.1048; "E36" = 50.2356; "F36" = 59.7389; "G36" = 132.6771; "H36" = 307.2284
}
foreach ($addr in $tbl4.Keys) {
  $ws.Range($addr).Value = $tbl4[$addr]
}

# A32:A36 carry the same bold-italic "n" column style as the other tables.
$ws.Range("A32:A36").Font.Bold = $true
$ws.Range("A32:A36").Font.Italic = $true

# ---------------------------------------------------------------------------
# 4) Conditional formatting: keep the two original 3-colour scales (now
#    lower priority) and add matching scales over the two newly-filled
#    tables (B23:H27 and B32:H36), highest priority first.
# ---------------------------------------------------------------------------
$ws.Range("B3:H7").FormatConditions.Item(1).Priority = 4
$ws.Range("B13:H17").FormatConditions.Item(1).Priority = 3

$csNew1 = $ws.Range("B32:H36").FormatConditions.AddColorScale(3)
$csNew1.ColorScaleCriteria(1).FormatColor.Color = 8109667
$csNew1.ColorScaleCriteria(2).FormatColor.Color = 8711167
$csNew1.ColorScaleCriteria(3).FormatColor.Color = 7039480
$csNew1.Priority = 2

$csNew2 = $ws.Range("B23:H27").FormatConditions.AddColorScale(3)
$csNew2.ColorScaleCriteria(1).FormatColor.Color = 8109667
$csNew2.ColorScaleCriteria(2).FormatColor.Color = 8711167
$csNew2.ColorScaleCriteria(3).FormatColor.Color = 7039480
$csNew2.Priority = 1

# ---------------------------------------------------------------------------
# 5) View state: zoom and final selection.
# ---------------------------------------------------------------------------
$excel.Windows.Item(1).Zoom = 118
$ws.Range("L31").Select() | Out-Null
